$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared date/time string used by column D (rows 3,4,5,7,10,13)
# from 2024-08-26T17:26:00.000Z to 2024-08-27T12:18:00.000Z
$ws.Range("D3").Value = "2024-08-27T12:18:00.000Z"
$ws.Range("D4").Value = "2024-08-27T12:18:00.000Z"
$ws.Range("D5").Value = "2024-08-27T12:18:00.000Z"
$ws.Range("D7").Value = "2024-08-27T12:18:00.000Z"
$ws.Range("D10").Value = "2024-08-27T12:18:00.000Z"
$ws.Range("D13").Value = "2024-08-27T12:18:00.000Z"

# Update numeric values on row 10 (chiet khau / thu no fix)
$ws.Range("T10").Value = 62500000
$ws.Range("W10").Value = 30575000
$ws.Range("AA10").Value = 89925000
$ws.Range("AE10").Value = 120500000
$ws.Range("AH10").Value = 98500000
$ws.Range("AK10").Value = 15
$ws.Range("AQ10").Value = 161000000
